# Quarterly indexing esoteric bug-fix operation
#
# The forecast-error table was re-indexed by one quarter: every existing
# data row (originally rows 2-10, corresponding to horizons Q0-Q8) shifted
# down into rows 3-11 (horizons Q1-Q9), while a brand new row of figures
# was computed for the first horizon and written into row 2 (Q0). The old
# last row (previously row 11 / Q9) is discarded since the table stays
# the same size.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for each data row, columns B:G (ME, MAE, MSE, RMSE, SE, N)
$data = @{
    2  = @(-0.02907897629796788, 0.3131278957257717, 0.181524606355785, 0.4260570458938391, 0.43998257208981, 15)
    3  = @(0.0370360702193368, 0.2877205798324301, 0.135345957066826, 0.3678939481247633, 0.3798420736582994, 14)
    4  = @(0.01989109342689383, 0.3298745939743749, 0.1503185731831121, 0.3877093926939507, 0.4030092977225764, 13)
    5  = @(0.06273519545459039, 0.2460334154802962, 0.09558537405740077, 0.3091688439306276, 0.3161984462764572, 12)
    6  = @(0.02711121047383215, 0.2779584568479502, 0.1216586105046287, 0.34879594393374, 0.3647135205176791, 11)
    7  = @(-0.01088135635153479, 0.2849421119723689, 0.09825224066429059, 0.3134521345664926, 0.3302084135617004, 10)
    8  = @(-0.04428645741563344, 0.3765628100937468, 0.1690784592812682, 0.4111915116843588, 0.4335975266214011, 9)
    9  = @(-0.01393931246739222, 0.3509192590318558, 0.1607817096315818, 0.4009759464501354, 0.4284022211487281, 8)
    10 = @(-0.03688841855209302, 0.2799590153621541, 0.1329006694222477, 0.3645554408073589, 0.3917438359423185, 7)
    11 = @(-0.06364682135181432, 0.1670219060428917, 0.03807752624502202, 0.1951346362002964, 0.2020691219662191, 6)
}

$columns = @("B", "C", "D", "E", "F", "G")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellRef = "$($columns[$i])$($row)"
        $ws.Range($cellRef).Value = $values[$i]
    }
}
